$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.276.98'
$ws.Range('E2').Value = '  +1.98%  '

$ws.Range('D3').Value = '2.097.57'
$ws.Range('E3').Value = '  +0.01%  '

$c = $ws.Range('D4')
$c.NumberFormat = '@'
$c.Value = '1.002'
$c.Style = 'Normal'
$ws.Range('E4').Value = '  -0.73%  '

$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '342.56'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -0.31%  '

$ws.Range('E6').Value = '  -0.62%  '

$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '0.5277'
$c.Style = 'Normal'
$ws.Range('E7').Value = '  +2.29%  '

$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.4381'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  +0.11%  '

$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '55.07'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  +2.99%  '

$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.09385'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  +2.31%  '

$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '1.177'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  +0.84%  '

$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '24.77'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  +0.90%  '

$ws.Range('E13').Value = '  +4.71%  '

$ws.Range('D14').Value = '2.113.75'
$ws.Range('E14').Value = '  +0.88%  '

$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '6.865'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  +1.52%  '

$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '101.31'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  -1.21%  '

$ws.Range('E17').Value = '  +0.41%  '

$ws.Range('E18').Value = '  -0.63%  '

$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '21.10'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  +0.57%  '

$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '0.06725'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  +0.81%  '

$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '6.414'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  +3.55%  '

$ws.Range('E22').Value = '  -0.58%  '

$ws.Range('D23').Value = '30.279.61'
$ws.Range('E23').Value = '  +1.74%  '

$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '12.45'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  -0.69%  '

$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '2.320'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  +0.67%  '

$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '7.013'
$c.Style = 'Normal'

$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '21.80'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  -0.54%  '

$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '162.49'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  +0.30%  '

$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '2.519'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  +0.84%  '

$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '133.77'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  +0.29%  '

$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '1.134'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  +0.40%  '

$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '1.676'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  +0.61%  '

$ws.Range('E33').Value = '  +0.24%  '

$ws.Range('E34').Value = '  +1.10%  '

$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '3.874'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  -2.27%  '

$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '0.02628'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  +2.15%  '

$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '0.06769'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  +1.20%  '

$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '12.69'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  +2.02%  '

$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '1.351'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  +1.93%  '

$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '0.6964'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  -0.46%  '

$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '0.2213'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  -0.14%  '

$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '0.6793'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  -0.30%  '

$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '14.33'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  -0.36%  '

$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '2.330'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  +0.73%  '

$ws.Range('E46').Value = '  -0.51%  '

$ws.Range('E47').Value = '  +8.51%  '

$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '3.632'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  +0.45%  '

$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '0.00000000346'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  -3.48%  '

$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '1.213'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  +6.57%  '

$ws.Range('E51').Value = '  +3.45%  '
